# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows at the top of the Kiwi "Macroferia Regional
# de Talca" block (rows 340-342), pushing the existing historical rows down
# by 3 (340-391 -> 343-394). Dimension grows from A1:T391 to A1:T394.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 340; Excel shifts every
# row at/after 340 down by three and copies formatting from the row above,
# so D341/D342 inherit the date style (s="2") automatically.
$ws.Rows("340:342").Insert()

function SetRow($r, $c1, $c2, $c3, $c4, $c5, $c6, $c7, $c8, $c9, $c10, $c11, $c12, $c13, $c14, $c15, $c16, $c17, $c18, $c19, $c20) {
    $ws.Cells.Item($r, 1).Value2 = $c1
    $ws.Cells.Item($r, 2).Value2 = $c2
    $ws.Cells.Item($r, 3).Value2 = $c3
    $ws.Cells.Item($r, 4).Value2 = $c4
    $ws.Cells.Item($r, 5).Value2 = $c5
    $ws.Cells.Item($r, 6).Value2 = $c6
    $ws.Cells.Item($r, 7).Value2 = $c7
    $ws.Cells.Item($r, 8).Value2 = $c8
    $ws.Cells.Item($r, 9).Value2 = $c9
    $ws.Cells.Item($r, 10).Value2 = $c10
    $ws.Cells.Item($r, 11).Value2 = $c11
    $ws.Cells.Item($r, 12).Value2 = $c12
    $ws.Cells.Item($r, 13).Value2 = $c13
    $ws.Cells.Item($r, 14).Value2 = $c14
    $ws.Cells.Item($r, 15).Value2 = $c15
    $ws.Cells.Item($r, 16).Value2 = $c16
    $ws.Cells.Item($r, 17).Value2 = $c17
    $ws.Cells.Item($r, 18).Value2 = $c18
    $ws.Cells.Item($r, 19).Value2 = $c19
    $ws.Cells.Item($r, 20).Value2 = $c20
}

SetRow 340 5 "Macroferia Regional de Talca" "Maule" 45034 7 "Fruta" 100101 "Berries" 100101007 "Kiwi" "Hayward" "Especial" 120 12000 12000 12000 "`$/bandeja 18 kilos" "Provincia de Curicó" 667 18

SetRow 341 5 "Macroferia Regional de Talca" "Maule" 45034 7 "Fruta" 100101 "Berries" 100101007 "Kiwi" "Hayward" "Primera" 260 10000 10000 10000 "`$/bandeja 18 kilos" "Provincia de Curicó" 556 18

SetRow 342 5 "Macroferia Regional de Talca" "Maule" 45034 7 "Fruta" 100101 "Berries" 100101007 "Kiwi" "Hayward" "Segunda" 250 8000 8000 8000 "`$/bandeja 18 kilos" "Provincia de Curicó" 444 18
